$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.003.59"
$ws.Range("D3").Value = "2.401.62"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "2.411.39"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "2.826.26"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "56.947.49"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.79"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "2.405.76"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "0.0₃0726"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.837"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "251.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0211"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.06"
$ws.Range("D51").ClearFormats()
